# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to the Leve profit sheets
# (computed columns H..N: currentAveragePrice*, LevePrice*, LeveProfit*)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 33666.332
$ws.Range("J16").Value = 33666.332
$ws.Range("L16").Value = 33666.332
$ws.Range("N16").Value = -34126.332
$ws.Range("H32").Value = 37038560
$ws.Range("I32").Value = 83333760
$ws.Range("J32").Value = 2400.4
$ws.Range("K32").Value = 83333760
$ws.Range("L32").Value = 2400.4
$ws.Range("M32").Value = -83333434
$ws.Range("N32").Value = -3052.4
$ws.Range("H132").Value = 23708.691
$ws.Range("I132").Value = 3361.697
$ws.Range("J132").Value = 135617.17
$ws.Range("K132").Value = 10085.091
$ws.Range("L132").Value = 406851.51
$ws.Range("M132").Value = -7555.091
$ws.Range("N132").Value = -411911.51
$ws.Range("H137").Value = 1976190.2
$ws.Range("I137").Value = 5495403.5
$ws.Range("J137").Value = 5430.88
$ws.Range("K137").Value = 16486210.5
$ws.Range("L137").Value = 16292.64
$ws.Range("M137").Value = -16483660.5
$ws.Range("N137").Value = -21392.64

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13562.022
$ws.Range("I32").Value = 12365.359
$ws.Range("J32").Value = 21340.334
$ws.Range("K32").Value = 12365.359
$ws.Range("L32").Value = 21340.334
$ws.Range("M32").Value = -12078.359
$ws.Range("N32").Value = -21914.334
$ws.Range("H102").Value = 40018336
$ws.Range("I102").Value = 47620644
$ws.Range("J102").Value = 106232.75
$ws.Range("K102").Value = 47620644
$ws.Range("L102").Value = 106232.75
$ws.Range("M102").Value = -47619022
$ws.Range("N102").Value = -109476.75
$ws.Range("H132").Value = 9261222
$ws.Range("I132").Value = 11629326
$ws.Range("J132").Value = 4086
$ws.Range("K132").Value = 34887978
$ws.Range("L132").Value = 12258
$ws.Range("M132").Value = -34885448
$ws.Range("N132").Value = -17318

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3101.1428
$ws.Range("I86").Value = 2901.889
$ws.Range("J86").Value = 3459.8
$ws.Range("K86").Value = 2901.889
$ws.Range("L86").Value = 3459.8
$ws.Range("M86").Value = -1778.889
$ws.Range("N86").Value = -5705.8
$ws.Range("H89").Value = 3101.1428
$ws.Range("I89").Value = 2901.889
$ws.Range("J89").Value = 3459.8
$ws.Range("K89").Value = 14509.445
$ws.Range("L89").Value = 17299
$ws.Range("M89").Value = -8893.445
$ws.Range("N89").Value = -28531

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4634314
$ws.Range("I31").Value = 1998.7858
$ws.Range("J31").Value = 7582151.5
$ws.Range("K31").Value = 1998.7858
$ws.Range("L31").Value = 7582151.5
$ws.Range("M31").Value = -1703.7858
$ws.Range("N31").Value = -7582741.5
$ws.Range("H34").Value = 4634314
$ws.Range("I34").Value = 1998.7858
$ws.Range("J34").Value = 7582151.5
$ws.Range("K34").Value = 1998.7858
$ws.Range("L34").Value = 7582151.5
$ws.Range("M34").Value = -1796.7858
$ws.Range("N34").Value = -7582555.5
$ws.Range("H62").Value = 2988.25
$ws.Range("I62").Value = 2850
$ws.Range("J62").Value = 3126.5
$ws.Range("K62").Value = 2850
$ws.Range("L62").Value = 3126.5
$ws.Range("M62").Value = -2226
$ws.Range("N62").Value = -4374.5
$ws.Range("H65").Value = 2988.25
$ws.Range("I65").Value = 2850
$ws.Range("J65").Value = 3126.5
$ws.Range("K65").Value = 14250
$ws.Range("L65").Value = 15632.5
$ws.Range("M65").Value = -11130
$ws.Range("N65").Value = -21872.5
$ws.Range("H132").Value = 33102.043
$ws.Range("I132").Value = 1086.9688
$ws.Range("J132").Value = 111908.38
$ws.Range("K132").Value = 3260.9064
$ws.Range("L132").Value = 335725.14
$ws.Range("M132").Value = -730.9064000000003
$ws.Range("N132").Value = -340785.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 280
$ws.Range("I36").Value = 280
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 840
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -671
$ws.Range("N36").ClearContents()
$ws.Range("H39").Value = 3372
$ws.Range("J39").Value = 3372
$ws.Range("L39").Value = 10116
$ws.Range("N39").Value = -10704
$ws.Range("H42").Value = 100
$ws.Range("I42").Value = 100
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 300
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 234
$ws.Range("N42").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H98").Value = 552.6
$ws.Range("J98").Value = 563.5
$ws.Range("L98").Value = 1690.5
$ws.Range("N98").Value = -4686.5
$ws.Range("H113").Value = 4100.6895
$ws.Range("I113").Value = 5735.5264
$ws.Range("J113").Value = 994.5
$ws.Range("K113").Value = 17206.5792
$ws.Range("L113").Value = 2983.5
$ws.Range("M113").Value = -15036.5792
$ws.Range("N113").Value = -7323.5
$ws.Range("H125").Value = 4407
$ws.Range("J125").Value = 3898.5454
$ws.Range("L125").Value = 11695.6362
$ws.Range("N125").Value = -21535.6362
$ws.Range("H131").Value = 812.7938
$ws.Range("I131").Value = 433.625
$ws.Range("J131").Value = 846.8764
$ws.Range("K131").Value = 1300.875
$ws.Range("L131").Value = 2540.6292
$ws.Range("M131").Value = 3739.125
$ws.Range("N131").Value = -12620.6292

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H70").Value = 4883.0625
$ws.Range("I70").Value = 4916.393
$ws.Range("J70").Value = 4649.75
$ws.Range("K70").Value = 4916.393
$ws.Range("L70").Value = 4649.75
$ws.Range("M70").Value = -4646.393
$ws.Range("N70").Value = -5189.75
$ws.Range("H73").Value = 4883.0625
$ws.Range("I73").Value = 4916.393
$ws.Range("J73").Value = 4649.75
$ws.Range("K73").Value = 4916.393
$ws.Range("L73").Value = 4649.75
$ws.Range("M73").Value = -3980.393
$ws.Range("N73").Value = -6521.75
$ws.Range("H132").Value = 43480916
$ws.Range("I132").Value = 71430190
$ws.Range("J132").Value = 4267
$ws.Range("K132").Value = 214290570
$ws.Range("L132").Value = 12801
$ws.Range("M132").Value = -214288040
$ws.Range("N132").Value = -17861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3088.3125
$ws.Range("I7").Value = 2168
$ws.Range("J7").Value = 3640.5
$ws.Range("K7").Value = 2168
$ws.Range("L7").Value = 3640.5
$ws.Range("M7").Value = -2056
$ws.Range("N7").Value = -3864.5
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H22").Value = 1110
$ws.Range("I22").Value = 1062.5
$ws.Range("K22").Value = 1062.5
$ws.Range("M22").Value = -767.5
$ws.Range("H25").Value = 10000
$ws.Range("J25").Value = 10000
$ws.Range("L25").Value = 10000
$ws.Range("N25").Value = -10460
$ws.Range("H27").Value = 1110
$ws.Range("I27").Value = 1062.5
$ws.Range("K27").Value = 1062.5
$ws.Range("M27").Value = -955.5
$ws.Range("H46").Value = 9088.388999999999
$ws.Range("I46").Value = 2350
$ws.Range("J46").Value = 11013.643
$ws.Range("K46").Value = 2350
$ws.Range("L46").Value = 11013.643
$ws.Range("M46").Value = -2162
$ws.Range("N46").Value = -11389.643
$ws.Range("H61").Value = 3050
$ws.Range("I61").Value = 3050
$ws.Range("K61").Value = 3050
$ws.Range("M61").Value = -2848
$ws.Range("H113").Value = 3050
$ws.Range("I113").Value = 3050
$ws.Range("K113").Value = 3050
$ws.Range("M113").Value = -880
$ws.Range("H126").Value = 3088.3125
$ws.Range("I126").Value = 2168
$ws.Range("J126").Value = 3640.5
$ws.Range("K126").Value = 6504
$ws.Range("L126").Value = 10921.5
$ws.Range("M126").Value = -4034
$ws.Range("N126").Value = -15861.5
$ws.Range("H132").Value = 3230.68
$ws.Range("I132").Value = 2306.3333
$ws.Range("K132").Value = 6918.999899999999
$ws.Range("M132").Value = -4388.999899999999
$ws.Range("H136").Value = 2151.5715
$ws.Range("I136").Value = 1721.0952
$ws.Range("J136").Value = 3443
$ws.Range("K136").Value = 5163.2856
$ws.Range("L136").Value = 10329
$ws.Range("M136").Value = -2613.2856
$ws.Range("N136").Value = -15429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 70005
$ws.Range("J7").Value = 70005
$ws.Range("L7").Value = 70005
$ws.Range("N7").Value = -70231
$ws.Range("H14").Value = 1000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 1000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 1000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -1336
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H19").Value = 37503
$ws.Range("I19").Value = 5000
$ws.Range("J19").Value = 70006
$ws.Range("K19").Value = 5000
$ws.Range("L19").Value = 70006
$ws.Range("M19").Value = -4826
$ws.Range("N19").Value = -70354
$ws.Range("H24").Value = 31501.5
$ws.Range("J24").Value = 34000.9
$ws.Range("L24").Value = 34000.9
$ws.Range("N24").Value = -34460.9
$ws.Range("H40").Value = 20827.143
$ws.Range("J40").Value = 20827.143
$ws.Range("L40").Value = 20827.143
$ws.Range("N40").Value = -21125.143
$ws.Range("H51").Value = 5713.3335
$ws.Range("I51").Value = 3570
$ws.Range("K51").Value = 3570
$ws.Range("M51").Value = -3060
$ws.Range("H69").Value = 17344.2
$ws.Range("J69").Value = 17344.2
$ws.Range("L69").Value = 17344.2
$ws.Range("N69").Value = -18842.2
$ws.Range("H72").Value = 17344.2
$ws.Range("J72").Value = 17344.2
$ws.Range("L72").Value = 52032.60000000001
$ws.Range("N72").Value = -59520.60000000001
$ws.Range("H132").Value = 1450967.9
$ws.Range("I132").Value = 2289312.8
$ws.Range("J132").Value = 2917.818
$ws.Range("K132").Value = 6867938.399999999
$ws.Range("L132").Value = 8753.454000000002
$ws.Range("M132").Value = -6865408.399999999
$ws.Range("N132").Value = -13813.454
